$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F39").Value = 198
$ws.Range("G39").Value = 5082.66
$ws.Range("F54").Value = 64
$ws.Range("G54").Value = 1966.08
$ws.Range("F55").Value = 37
$ws.Range("G55").Value = 681.54
$ws.Range("F68").Value = 49
$ws.Range("G68").Value = 1730.68
$ws.Range("B71").Value = 61436.23
$ws.Range("F141").Value = 491
$ws.Range("G141").Value = 9554.860000000001
$ws.Range("B143").Value = 295897.16
$ws.Range("F167").Value = 64
$ws.Range("G167").Value = 6665.6
$ws.Range("B176").Value = 17847.9
$ws.Range("F196").Value = 8
$ws.Range("G196").Value = 338.56
$ws.Range("B205").Value = 28715.99
$ws.Range("F233").Value = 5
$ws.Range("G233").Value = 1342.5
$ws.Range("F234").Value = 34
$ws.Range("G234").Value = 2839
$ws.Range("F237").Value = 18
$ws.Range("G237").Value = 2003.94
$ws.Range("F242").Value = 22
$ws.Range("G242").Value = 2655.18
$ws.Range("F249").Value = 16
$ws.Range("G249").Value = 2052.48
$ws.Range("B251").Value = 36685.32
$ws.Range("B388").Value = 61610
$ws.Range("D388").Value = 102.71
$ws.Range("E388").Value = 122.71
$ws.Range("F388").Value = 266
$ws.Range("G388").Value = 27320.86
$ws.Range("B389").Value = 57077
$ws.Range("D389").Value = 93.08
$ws.Range("E389").Value = 111.2
$ws.Range("F389").Value = 1
$ws.Range("G389").Value = 93.08
$ws.Range("F392").Value = 91
$ws.Range("G392").Value = 9745.190000000001
$ws.Range("F398").Value = 30
$ws.Range("G398").Value = 3619.5
$ws.Range("F422").Value = 193
$ws.Range("G422").Value = 21608.28
$ws.Range("F423").Value = 190
$ws.Range("G423").Value = 24511.9
$ws.Range("F439").Value = 13
$ws.Range("G439").Value = 708.24
$ws.Range("F441").Value = 16
$ws.Range("G441").Value = 1662.56
$ws.Range("F454").Value = 40
$ws.Range("G454").Value = 7439.6
$ws.Range("F455").Value = 278
$ws.Range("G455").Value = 19546.18
$ws.Range("F457").Value = 56
$ws.Range("G457").Value = 8465.52
$ws.Range("B468").Value = 458015.41
$ws.Range("F566").Value = 73
$ws.Range("G566").Value = 4412.85
$ws.Range("B572").Value = 60401.23
$ws.Range("F590").Value = 652
$ws.Range("G590").Value = 8352.120000000001
$ws.Range("F595").Value = 402
$ws.Range("G595").Value = 7931.46
$ws.Range("F602").Value = 454
$ws.Range("G602").Value = 5970.1
$ws.Range("F605").Value = 654
$ws.Range("G605").Value = 9633.42
$ws.Range("B606").Value = 114422.29
$ws.Range("F617").Value = 10
$ws.Range("G617").Value = 348.1
$ws.Range("B621").Value = 13788.18
$ws.Range("F672").Value = 388
$ws.Range("G672").Value = 5001.32
$ws.Range("F673").Value = 681
$ws.Range("G673").Value = 13517.85
$ws.Range("F674").Value = 480
$ws.Range("G674").Value = 3216
$ws.Range("F677").Value = 161
$ws.Range("G677").Value = 5255.04
$ws.Range("B678").Value = 44605.8
$ws.Range("F725").Value = 57
$ws.Range("G725").Value = 1260.27
$ws.Range("F738").Value = 74
$ws.Range("G738").Value = 321.9
$ws.Range("B740").Value = 10103.02
$ws.Range("F743").Value = 5
$ws.Range("G743").Value = 1429.35
$ws.Range("F756").Value = 4
$ws.Range("G756").Value = 533.3200000000001
$ws.Range("F759").Value = 91
$ws.Range("G759").Value = 11177.53
$ws.Range("B762").Value = 63753.74
$ws.Range("F797").Value = 322
$ws.Range("G797").Value = 5106.92
$ws.Range("F799").Value = 203
$ws.Range("G799").Value = 8765.540000000001
$ws.Range("F801").Value = 120
$ws.Range("G801").Value = 5181.6
$ws.Range("B805").Value = 41926.02
$ws.Range("F821").Value = 13
$ws.Range("G821").Value = 8851.959999999999
$ws.Range("B827").Value = 25072.32
$ws.Range("F835").Value = 30
$ws.Range("G835").Value = 2565
$ws.Range("B840").Value = 10451.65
$ws.Range("F864").Value = 28
$ws.Range("G864").Value = 2455.6
$ws.Range("B867").Value = 5801.27
$ws.Range("F894").Value = 112
$ws.Range("G894").Value = 4642.4
$ws.Range("F896").Value = 25
$ws.Range("G896").Value = 3778.25
$ws.Range("B902").Value = 69410.69
$ws.Range("F921").Value = 260
$ws.Range("G921").Value = 7831.2
$ws.Range("F926").Value = 19
$ws.Range("G926").Value = 1318.22
$ws.Range("F927").Value = 60
$ws.Range("G927").Value = 2209.8
$ws.Range("B934").Value = 48606.18
$ws.Range("F937").Value = 93
$ws.Range("G937").Value = 3478.2
$ws.Range("F940").Value = 174
$ws.Range("G940").Value = 6507.6
$ws.Range("F942").Value = 157
$ws.Range("G942").Value = 5871.8
$ws.Range("B943").Value = 18727.7
$ws.Range("F987").Value = 6
$ws.Range("G987").Value = 247.56
$ws.Range("B998").Value = 4782.86
$ws.Range("F1002").Value = 246
$ws.Range("G1002").Value = 18976.44
$ws.Range("F1005").Value = 16
$ws.Range("G1005").Value = 2057.28
$ws.Range("B1006").Value = 338540.2
$ws.Range("B1013").Value = 2963462.6
$ws.Range("B1014").Value = 2963462.6
